# "Final restructuring before golive"
# Adds three new articles (Bienvenido / Adopcion / encontre una mascota) at
# the top of the list, shifts the "Cuales medios de pago?" FAQ block down,
# bumps the `categories` id (B column) for every block, gives B2 its own
# highlight formatting, and normalizes a few cells that had stray styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row: C1 ("icon") had a stray/empty-font style; align it with the
# rest of the header row (A1/B1/D1/E1).
# ---------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# Row 2: new "Bienvenido" article
# ---------------------------------------------------------------------
$ws.Range("B2").Value = 21
$ws.Range("C2").Value = "bienvenido.jpg"
$ws.Range("D2").Value = "Bienvenido"
# E2 (body) is unchanged - still the Lorem ipsum placeholder.

# give B2 its own highlight style: 11pt black Arial on a solid white fill
$ws.Range("B2").Font.Size = 11
$ws.Range("B2").Font.Name = "Arial"
$ws.Range("B2").Font.ColorIndex = 1
$ws.Range("B2").Interior.Pattern = 1
$ws.Range("B2").Interior.Color = 16777215
$ws.Range("B2").Interior.PatternColor = 16777215

# match the new cells' look to their neighbours in the same column
$ws.Range("D2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# Row 3: new "Adopcion" article
# ---------------------------------------------------------------------
$ws.Range("B3:B4").Formula = "=B2"
$ws.Range("C3").Value = "adoptame.jpg"
$ws.Range("D3").Value = "Adopcion"

$ws.Range("D3").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# Row 4: new "encontre una mascota" article
# ---------------------------------------------------------------------
$ws.Range("C4").Value = "encontre_o_perdi_mascota.jpg"
$ws.Range("D4").Value = "encontre una mascota"

$ws.Range("D4").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# Row 5: first row of the (now-shifted) "Cuales medios de pago?" FAQ
# block. Breaks the old shared-formula chain in column B with a literal
# value, same as the source workbook's new layout.
# ---------------------------------------------------------------------
$ws.Range("B5").Value = 20
$ws.Range("D5").Value = "Cuales medios de pago?"

# ---------------------------------------------------------------------
# Rows 6-15: rest of the FAQ block - category id now fills down from B5.
# ---------------------------------------------------------------------
$ws.Range("B6:B15").Formula = "=B5"
$ws.Range("D6").Value = "Cuales medios de pago?"
$ws.Range("D7").Value = "Cuales medios de pago?"
$ws.Range("D8").Value = "Cuales medios de pago?"
$ws.Range("D9").Value = "Cuales medios de pago?"
$ws.Range("D10").Value = "Cuales medios de pago?"
$ws.Range("D11").Value = "Cuales medios de pago?"
$ws.Range("D12").Value = "Cuales medios de pago?"
$ws.Range("D13").Value = "Cuales medios de pago?"
$ws.Range("D14").Value = "Cuales medios de pago?"
$ws.Range("D15").Value = "Cuales medios de pago?"

# ---------------------------------------------------------------------
# Rows 16-18: Banner block - category id bumped 14 -> 21.
# ---------------------------------------------------------------------
$ws.Range("B16").Value = 21
$ws.Range("B17").Value = 21
$ws.Range("B18").Value = 21

# D17/D18 previously carried a stray style; normalize to the plain look
# used everywhere else in column D.
$ws.Range("D16").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D16").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# Rows 19-21: Noticia block - category id bumped 15 -> 22.
# ---------------------------------------------------------------------
$ws.Range("B19").Value = 22
$ws.Range("B20").Value = 22
$ws.Range("B21").Value = 22
